# Update the days vocabulary prefix URI from a "#" delimiter to a "/" delimiter,
# and move the active selection on "Feuil2" to C2 (matching the saved workbook
# view state after this edit).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Feuil2")
$ws.Activate()

$ws.Range("C2").Value = "http://data.sparna.fr/vocabularies/days/"

$ws.Range("C2").Select() | Out-Null
